# Script: applies the betexplorer laliga 2023-2024 refresh
#  1) swap match rows 58 <-> 59 (data columns F:V; A:E stay keyed to row position)
#  2) rotate match rows 63 -> 64 -> 65 -> 63 (data columns F:V)
#  3) swap match rows 68 <-> 69 (data columns F:V)
#  4) append a new match row 116 (Alaves vs Almeria) at the bottom of the sheet

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) swap rows 58 and 59 (columns F:V hold the per-match data) ---
$row58 = $ws.Range("F58:V58").Value()
$row59 = $ws.Range("F59:V59").Value()
$ws.Range("F58:V58").Value = $row59
$ws.Range("F59:V59").Value = $row58

# --- 2) rotate rows 63, 64, 65 ---
$row63 = $ws.Range("F63:V63").Value()
$row64 = $ws.Range("F64:V64").Value()
$row65 = $ws.Range("F65:V65").Value()
$ws.Range("F63:V63").Value = $row65
$ws.Range("F64:V64").Value = $row63
$ws.Range("F65:V65").Value = $row64

# --- 3) swap rows 68 and 69 ---
$row68 = $ws.Range("F68:V68").Value()
$row69 = $ws.Range("F69:V69").Value()
$ws.Range("F68:V68").Value = $row69
$ws.Range("F69:V69").Value = $row68

# --- 4) append new row 116 with the same look & feel as the existing rows ---
# Copy the number-formats from row 115 (index column bold/centered, date column
# date-time format) onto row 116 before writing the values.
$ws.Range("A115").Copy()
$ws.Range("A116").PasteSpecial(-4122)
$ws.Range("E115").Copy()
$ws.Range("E116").PasteSpecial(-4122)

$ws.Range("A116").Value = 115
$ws.Range("B116").Value = "spain"
$ws.Range("C116").Value = "laliga"
$ws.Range("D116").Value = "2023-2024"
$ws.Range("E116").Value = 45235.58333333334
$ws.Range("F116").Value = "Alaves"
$ws.Range("G116").Value = 1
$ws.Range("H116").Value = "Almeria"
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 1.87
$ws.Range("K116").Value = "23/10/2023 15:49"
$ws.Range("L116").Value = 1.77
$ws.Range("M116").Value = "05/11/2023 13:32"
$ws.Range("N116").Value = 3.64
$ws.Range("O116").Value = "23/10/2023 15:49"
$ws.Range("P116").Value = 3.77
$ws.Range("Q116").Value = "05/11/2023 13:32"
$ws.Range("R116").Value = 4.43
$ws.Range("S116").Value = "23/10/2023 15:49"
$ws.Range("T116").Value = 4.99
$ws.Range("U116").Value = "05/11/2023 13:32"
$ws.Range("V116").Value = "https://www.betexplorer.com/football/spain/laliga/alaves-almeria/pbJlLa6P/"

Write-Host "Applied laliga 2023-2024 update (rows 58/59 swap, 63-65 rotate, 68/69 swap, row 116 appended)"
